# Update NATMI LR-pair stats (Postn-Ptk7) with recomputed TPM-based values.
# Commit: "update scripts wuth new tpm"
# Applies the exact numeric changes to columns G:T for rows 2-17.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"20.66848433333334"
$ws.Range("H2").Value = [double]"62.005453"
$ws.Range("I2").Value = [double]"0.004921559247345127"
$ws.Range("J2").Value = [double]"0.004921559247345126"
$ws.Range("K2").Value = [double]"3"
$ws.Range("L2").Value = [double]"1"
$ws.Range("M2").Value = [double]"4.609402999999999"
$ws.Range("N2").Value = [double]"13.828209"
$ws.Range("O2").Value = [double]"0.1736766757621145"
$ws.Range("P2").Value = [double]"0.1736766757621145"
$ws.Range("Q2").Value = [double]"95.26937369151966"
$ws.Range("R2").Value = [double]"857.424363223677"
$ws.Range("S2").Value = [double]"0.0008547600496451959"
$ws.Range("T2").Value = [double]"0.0008547600496451958"
$ws.Range("G3").Value = [double]"20.66848433333334"
$ws.Range("H3").Value = [double]"62.005453"
$ws.Range("I3").Value = [double]"0.004921559247345127"
$ws.Range("J3").Value = [double]"0.004921559247345126"
$ws.Range("O3").Value = [double]"0.428411030010331"
$ws.Range("P3").Value = [double]"0.4284110300103309"
$ws.Range("Q3").Value = [double]"235.0024856966214"
$ws.Range("R3").Value = [double]"2115.022371269592"
$ws.Range("S3").Value = [double]"0.002108450266411995"
$ws.Range("T3").Value = [double]"0.002108450266411994"
$ws.Range("G4").Value = [double]"20.66848433333334"
$ws.Range("H4").Value = [double]"62.005453"
$ws.Range("I4").Value = [double]"0.004921559247345127"
$ws.Range("J4").Value = [double]"0.004921559247345126"
$ws.Range("M4").Value = [double]"10.551036"
$ws.Range("N4").Value = [double]"31.653108"
$ws.Range("O4").Value = [double]"0.3975501509254882"
$ws.Range("P4").Value = [double]"0.3975501509254881"
$ws.Range("Q4").Value = [double]"218.0739222664361"
$ws.Range("R4").Value = [double]"1962.665300397924"
$ws.Range("S4").Value = [double]"0.001956566621570787"
$ws.Range("T4").Value = [double]"0.001956566621570786"
$ws.Range("G5").Value = [double]"20.66848433333334"
$ws.Range("H5").Value = [double]"62.005453"
$ws.Range("I5").Value = [double]"0.004921559247345127"
$ws.Range("J5").Value = [double]"0.004921559247345126"
$ws.Range("K5").Value = [double]"1"
$ws.Range("L5").Value = [double]"0.3333333333333333"
$ws.Range("M5").Value = [double]"0.009611333333333333"
$ws.Range("N5").Value = [double]"0.028834"
$ws.Range("O5").Value = [double]"0.0003621433020664361"
$ws.Range("P5").Value = [double]"0.000362143302066436"
$ws.Range("Q5").Value = [double]"0.1986516924224445"
$ws.Range("R5").Value = [double]"1.787865231802"
$ws.Range("S5").Value = [double]"1.782309717149168E-06"
$ws.Range("T5").Value = [double]"1.782309717149168E-06"
$ws.Range("I6").Value = [double]"0.9485231866833483"
$ws.Range("J6").Value = [double]"0.9485231866833483"
$ws.Range("K6").Value = [double]"3"
$ws.Range("L6").Value = [double]"1"
$ws.Range("M6").Value = [double]"4.609402999999999"
$ws.Range("N6").Value = [double]"13.828209"
$ws.Range("O6").Value = [double]"0.1736766757621145"
$ws.Range("P6").Value = [double]"0.1736766757621145"
$ws.Range("Q6").Value = [double]"18361.09358552442"
$ws.Range("R6").Value = [double]"165249.8422697198"
$ws.Range("S6").Value = [double]"0.1647363539464515"
$ws.Range("T6").Value = [double]"0.1647363539464515"
$ws.Range("I7").Value = [double]"0.9485231866833483"
$ws.Range("J7").Value = [double]"0.9485231866833483"
$ws.Range("O7").Value = [double]"0.428411030010331"
$ws.Range("P7").Value = [double]"0.4284110300103309"
$ws.Range("S7").Value = [double]"0.4063577953956947"
$ws.Range("T7").Value = [double]"0.4063577953956946"
$ws.Range("I8").Value = [double]"0.9485231866833483"
$ws.Range("J8").Value = [double]"0.9485231866833483"
$ws.Range("M8").Value = [double]"10.551036"
$ws.Range("N8").Value = [double]"31.653108"
$ws.Range("O8").Value = [double]"0.3975501509254882"
$ws.Range("P8").Value = [double]"0.3975501509254881"
$ws.Range("Q8").Value = [double]"42028.99148116085"
$ws.Range("R8").Value = [double]"378260.9233304476"
$ws.Range("S8").Value = [double]"0.3770855360222902"
$ws.Range("T8").Value = [double]"0.3770855360222901"
$ws.Range("I9").Value = [double]"0.9485231866833483"
$ws.Range("J9").Value = [double]"0.9485231866833483"
$ws.Range("K9").Value = [double]"1"
$ws.Range("L9").Value = [double]"0.3333333333333333"
$ws.Range("M9").Value = [double]"0.009611333333333333"
$ws.Range("N9").Value = [double]"0.028834"
$ws.Range("O9").Value = [double]"0.0003621433020664361"
$ws.Range("P9").Value = [double]"0.000362143302066436"
$ws.Range("Q9").Value = [double]"38.28578035268422"
$ws.Range("R9").Value = [double]"344.5720231741579"
$ws.Range("S9").Value = [double]"0.0003435013189120864"
$ws.Range("T9").Value = [double]"0.0003435013189120863"
$ws.Range("G10").Value = [double]"192.4566396666667"
$ws.Range("H10").Value = [double]"577.369919"
$ws.Range("I10").Value = [double]"0.04582758655103054"
$ws.Range("J10").Value = [double]"0.04582758655103054"
$ws.Range("K10").Value = [double]"3"
$ws.Range("L10").Value = [double]"1"
$ws.Range("M10").Value = [double]"4.609402999999999"
$ws.Range("N10").Value = [double]"13.828209"
$ws.Range("O10").Value = [double]"0.1736766757621145"
$ws.Range("P10").Value = [double]"0.1736766757621145"
$ws.Range("Q10").Value = [double]"887.1102122494522"
$ws.Range("R10").Value = [double]"7983.99191024507"
$ws.Range("S10").Value = [double]"0.007959182890383571"
$ws.Range("T10").Value = [double]"0.007959182890383571"
$ws.Range("G11").Value = [double]"192.4566396666667"
$ws.Range("H11").Value = [double]"577.369919"
$ws.Range("I11").Value = [double]"0.04582758655103054"
$ws.Range("J11").Value = [double]"0.04582758655103054"
$ws.Range("O11").Value = [double]"0.428411030010331"
$ws.Range("P11").Value = [double]"0.4284110300103309"
$ws.Range("Q11").Value = [double]"2188.248929194291"
$ws.Range("R11").Value = [double]"19694.24036274862"
$ws.Range("S11").Value = [double]"0.01963304355721459"
$ws.Range("T11").Value = [double]"0.01963304355721458"
$ws.Range("G12").Value = [double]"192.4566396666667"
$ws.Range("H12").Value = [double]"577.369919"
$ws.Range("I12").Value = [double]"0.04582758655103054"
$ws.Range("J12").Value = [double]"0.04582758655103054"
$ws.Range("M12").Value = [double]"10.551036"
$ws.Range("N12").Value = [double]"31.653108"
$ws.Range("O12").Value = [double]"0.3975501509254882"
$ws.Range("P12").Value = [double]"0.3975501509254881"
$ws.Range("Q12").Value = [double]"2030.616933562028"
$ws.Range("R12").Value = [double]"18275.55240205825"
$ws.Range("S12").Value = [double]"0.01821876394991307"
$ws.Range("T12").Value = [double]"0.01821876394991306"
$ws.Range("G13").Value = [double]"192.4566396666667"
$ws.Range("H13").Value = [double]"577.369919"
$ws.Range("I13").Value = [double]"0.04582758655103054"
$ws.Range("J13").Value = [double]"0.04582758655103054"
$ws.Range("K13").Value = [double]"1"
$ws.Range("L13").Value = [double]"0.3333333333333333"
$ws.Range("M13").Value = [double]"0.009611333333333333"
$ws.Range("N13").Value = [double]"0.028834"
$ws.Range("O13").Value = [double]"0.0003621433020664361"
$ws.Range("P13").Value = [double]"0.000362143302066436"
$ws.Range("Q13").Value = [double]"1.849764916049556"
$ws.Range("R13").Value = [double]"16.647884244446"
$ws.Range("S13").Value = [double]"1.65961535193256E-05"
$ws.Range("T13").Value = [double]"1.659615351932559E-05"
$ws.Range("G14").Value = [double]"3.055898333333333"
$ws.Range("H14").Value = [double]"9.167694999999998"
$ws.Range("I14").Value = [double]"0.0007276675182760082"
$ws.Range("J14").Value = [double]"0.0007276675182760081"
$ws.Range("K14").Value = [double]"3"
$ws.Range("L14").Value = [double]"1"
$ws.Range("M14").Value = [double]"4.609402999999999"
$ws.Range("N14").Value = [double]"13.828209"
$ws.Range("O14").Value = [double]"0.1736766757621145"
$ws.Range("P14").Value = [double]"0.1736766757621145"
$ws.Range("Q14").Value = [double]"14.08586694536166"
$ws.Range("R14").Value = [double]"126.772802508255"
$ws.Range("S14").Value = [double]"0.0001263788756342448"
$ws.Range("T14").Value = [double]"0.0001263788756342448"
$ws.Range("G15").Value = [double]"3.055898333333333"
$ws.Range("H15").Value = [double]"9.167694999999998"
$ws.Range("I15").Value = [double]"0.0007276675182760082"
$ws.Range("J15").Value = [double]"0.0007276675182760081"
$ws.Range("O15").Value = [double]"0.428411030010331"
$ws.Range("P15").Value = [double]"0.4284110300103309"
$ws.Range("Q15").Value = [double]"34.74583296905333"
$ws.Range("R15").Value = [double]"312.71249672148"
$ws.Range("S15").Value = [double]"0.000311740791009686"
$ws.Range("T15").Value = [double]"0.0003117407910096859"
$ws.Range("G16").Value = [double]"3.055898333333333"
$ws.Range("H16").Value = [double]"9.167694999999998"
$ws.Range("I16").Value = [double]"0.0007276675182760082"
$ws.Range("J16").Value = [double]"0.0007276675182760081"
$ws.Range("M16").Value = [double]"10.551036"
$ws.Range("N16").Value = [double]"31.653108"
$ws.Range("O16").Value = [double]"0.3975501509254882"
$ws.Range("P16").Value = [double]"0.3975501509254881"
$ws.Range("Q16").Value = [double]"32.24289332734"
$ws.Range("R16").Value = [double]"290.18603994606"
$ws.Range("S16").Value = [double]"0.0002892843317142026"
$ws.Range("T16").Value = [double]"0.0002892843317142024"
$ws.Range("G17").Value = [double]"3.055898333333333"
$ws.Range("H17").Value = [double]"9.167694999999998"
$ws.Range("I17").Value = [double]"0.0007276675182760082"
$ws.Range("J17").Value = [double]"0.0007276675182760081"
$ws.Range("K17").Value = [double]"1"
$ws.Range("L17").Value = [double]"0.3333333333333333"
$ws.Range("M17").Value = [double]"0.009611333333333333"
$ws.Range("N17").Value = [double]"0.028834"
$ws.Range("O17").Value = [double]"0.0003621433020664361"
$ws.Range("P17").Value = [double]"0.000362143302066436"
$ws.Range("Q17").Value = [double]"0.02937125751444444"
$ws.Range("R17").Value = [double]"0.26434131763"
$ws.Range("S17").Value = [double]"2.635199178749624E-07"
$ws.Range("T17").Value = [double]"2.635199178749623E-07"
